# Update the lattice-multiplication exercise table: every cell in the
# 5-row x 3-column table gets a new multiplication problem (same 5-line
# lattice layout: "A x B" / "  B    B" factor row / "  ----" separator /
# two partial-product placeholder rows).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$br = [char]11

$t.Cell(1, 1).Range.Text = "36 x 74" + $br + "  7    4" + $br + "  ----" + $br + "3|    |" + $br + "6|    |"
$t.Cell(1, 2).Range.Text = "31 x 35" + $br + "  3    5" + $br + "  ----" + $br + "3|    |" + $br + "1|    |"
$t.Cell(1, 3).Range.Text = "31 x 51" + $br + "  5    1" + $br + "  ----" + $br + "3|    |" + $br + "1|    |"

$t.Cell(2, 1).Range.Text = "23 x 60" + $br + "  6    0" + $br + "  ----" + $br + "2|    |" + $br + "3|    |"
$t.Cell(2, 2).Range.Text = "60 x 23" + $br + "  2    3" + $br + "  ----" + $br + "6|    |" + $br + "0|    |"
$t.Cell(2, 3).Range.Text = "25 x 44" + $br + "  4    4" + $br + "  ----" + $br + "2|    |" + $br + "5|    |"

$t.Cell(3, 1).Range.Text = "30 x 53" + $br + "  5    3" + $br + "  ----" + $br + "3|    |" + $br + "0|    |"
$t.Cell(3, 2).Range.Text = "88 x 87" + $br + "  8    7" + $br + "  ----" + $br + "8|    |" + $br + "8|    |"
$t.Cell(3, 3).Range.Text = "48 x 87" + $br + "  8    7" + $br + "  ----" + $br + "4|    |" + $br + "8|    |"

$t.Cell(4, 1).Range.Text = "72 x 12" + $br + "  1    2" + $br + "  ----" + $br + "7|    |" + $br + "2|    |"
$t.Cell(4, 2).Range.Text = "75 x 49" + $br + "  4    9" + $br + "  ----" + $br + "7|    |" + $br + "5|    |"
$t.Cell(4, 3).Range.Text = "21 x 63" + $br + "  6    3" + $br + "  ----" + $br + "2|    |" + $br + "1|    |"

$t.Cell(5, 1).Range.Text = "16 x 10" + $br + "  1    0" + $br + "  ----" + $br + "1|    |" + $br + "6|    |"
$t.Cell(5, 2).Range.Text = "91 x 26" + $br + "  2    6" + $br + "  ----" + $br + "9|    |" + $br + "1|    |"
$t.Cell(5, 3).Range.Text = "74 x 88" + $br + "  8    8" + $br + "  ----" + $br + "7|    |" + $br + "4|    |"

Write-Host "Lattice table cells updated."
